$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.700.08'
$ws.Range('D3').Value = '1.592.59'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '208.76'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D8').Value = '22.30'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').Value = '0.0868'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '1.818.12'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').Value = '1.579.07'
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').Value = '0.529'
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '27.694.46'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').Value = '63.27'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '218.02'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '4.16'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').Value = '9.79'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').Value = '153.87'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D28').Value = '15.09'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = '0.0475'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('D33').Value = '1.382.23'
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = '2.32'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('D39').Value = '0.535'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '0.828'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('D43').Value = '64.46'
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('E44').Value = '  +4.40%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '1.76'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '5.26'
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('D47').Value = '1.729.66'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').Value = '85.86'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('E51').Value = '  -0.19%  '
